$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" header column in H1, matching the formatting used by
# the other header cells (bold, bordered, centered) by copying the
# format from the existing "sum" header (G1).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add the corresponding value for row 2 (Save = 0)
$ws.Range("H2").Value = 0
